# ============================================================================
# Applies the "acervo_6-2" / "acervo_6-3" glossary-sheet addition described
# by the commit "Atualizações até a data de 13/01/2022."
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the two new worksheets, right after acervo_6-1 ---------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "acervo_6-2"
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "acervo_6-3"

# --- acervo_6-2 : variable glossary ---------------------------------------
$ws2.Range("A1").Value = "Variável"
$ws2.Range("B1").Value = "Definição"
$ws2.Range("A2").Value = "Id Processo"
$ws2.Range("B2").Value = "Chave numérica única para cada processo."
$ws2.Range("A3").Value = "Classe"
$ws2.Range("B3").Value = "Sigla da classe do processo no STF."
$ws2.Range("A4").Value = "Número"
$ws2.Range("B4").Value = "Número do processo no STF."
$ws2.Range("A5").Value = "Nome Ministro (a)"
$ws2.Range("B5").Value = "Informa para qual Ministro o processo foi distribuído"
$ws2.Range("A6").Value = "Data Andamento"
$ws2.Range("B6").Value = "Data em que o processo foi distribuído"
$ws2.Range("A7").Value = "Qtd Ocorrências Processuais"
$ws2.Range("B7").Value = "Informa a quantidade de ocorrências em cada andamento processual."
$ws2.Range("A8").Value = "Andamento"
$ws2.Range("B8").Value = "Descreve qual andamento foi lançado no ato da distribuição."
$ws2.Range("A9").Value = "Classificação STF"
$ws2.Range("B9").Value = "Recursal ou originária"
$ws2.Range("A10").Value = "Subgrupo Andamento"
$ws2.Range("B10").Value = "Informa a qual Subgrupo o andamento pertence."

# --- acervo_6-3 : chart/table glossary ------------------------------------
$ws3.Range("A1").Value = "Nome do gráfico e/ou tabela"
$ws3.Range("B1").Value = "Definição"
$ws3.Range("A2").Value = "Registrados e Distribuídos"
$ws3.Range("B2").Value = "Quantitativos geral dos processos registrados à Presidência e distribuídos aos Ministros"
$ws3.Range("A3").Value = "Registro à Presidência"
$ws3.Range("B3").Value = "Quantitativo dos processos registrados à Presidência"
$ws3.Range("A4").Value = "Distribuídos aos Ministros"
$ws3.Range("B4").Value = "Quantitativo dos processos distribuídos aos Ministros"
$ws3.Range("A5").Value = "Distribuídos e Registrados à Presidência"
$ws3.Range("B5").Value = "Gráfico em linha dos processos distribuídos por ano"
$ws3.Range("A6").Value = "Processos pro classe"
$ws3.Range("B6").Value = "Gráfico em barra dos processos por classe processual"
$ws3.Range("A7").Value = "Processos Distribuídos aos Ministros"
$ws3.Range("B7").Value = "Gráfico em barra da distribuição por ministro"
$ws3.Range("A8").Value = "Gráfico pizza"
$ws3.Range("B8").Value = "Distribuição por processo recursal ou originário"
$ws3.Range("A9").Value = "Processos por ramo do direito"
$ws3.Range("B9").Value = "Mapa de árvore da distribuição por ramo do direito"
$ws3.Range("A10").Value = "Processos Distribuídos por órgão origem"
$ws3.Range("B10").Value = "Tabela dos processos distribuídos aos Ministros por órgão origem"
$ws3.Range("A11").Value = "Processo Registrados por órgão origem"
$ws3.Range("B11").Value = "Tabela dos processos registrados à Presidência por órgão origem"

# --- formatting: reuse the same "label" look used on acervo_6-1 ----------
foreach ($ws in @($ws2, $ws3)) {
    $used = $ws.UsedRange
    $used.Font.Name = "Segoe UI"
    $used.Font.Size = 8
    $used.Font.Color = 5129531
    $used.HorizontalAlignment = -4131
    $used.VerticalAlignment = -4108

    $ws.PageSetup.PaperSize = 9
    $ws.PageSetup.Orientation = 1
    $ws.PageSetup.LeftMargin = 36.850393728
    $ws.PageSetup.RightMargin = 36.850393728
    $ws.PageSetup.TopMargin = 56.692913399999995
    $ws.PageSetup.BottomMargin = 56.692913399999995
    $ws.PageSetup.HeaderMargin = 22.67716464
    $ws.PageSetup.FooterMargin = 22.67716464
}

# column widths (character units, minus the ~5px/6 padding baked into the
# stored XML "width" so the saved file lands on the target values)
$ws2.Columns.Item(1).ColumnWidth = 18.166666666666668
$ws2.Columns.Item(2).ColumnWidth = 44.529947916666664
$ws3.Columns.Item(1).ColumnWidth = 29.436197916666668
$ws3.Columns.Item(2).ColumnWidth = 63.529947916666664

# --- sheet views: zoom + selection -----------------------------------------
$ws2.Activate()
$excel.ActiveWindow.Zoom = 115
$ws2.Range("E9").Select() | Out-Null

$ws3.Activate()
$ws3.Range("E7").Select() | Out-Null

Write-Host "done"
